$wb = $excel.ActiveWorkbook

# ---- ALC (sheet1.xml) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5289
$ws.Range("J17").Value = 5281.531
$ws.Range("L17").Value = 15844.593
$ws.Range("N17").Value = -16180.593
$ws.Range("H28").Value = 1111.8158
$ws.Range("I28").Value = 378.51852
$ws.Range("J28").Value = 2911.7273
$ws.Range("K28").Value = 378.51852
$ws.Range("L28").Value = 2911.7273
$ws.Range("M28").Value = 106.48148
$ws.Range("N28").Value = -3881.7273
$ws.Range("H40").Value = 11476.667
$ws.Range("J40").Value = 5153.5713
$ws.Range("L40").Value = 5153.5713
$ws.Range("N40").Value = -5503.5713
$ws.Range("H53").Value = 12974.5625
$ws.Range("I53").Value = 455.85715
$ws.Range("K53").Value = 455.85715
$ws.Range("M53").Value = 181.14285
$ws.Range("H61").Value = 282.4
$ws.Range("I61").Value = 282.4
$ws.Range("K61").Value = 847.1999999999999
$ws.Range("M61").Value = -675.1999999999999
$ws.Range("H64").Value = 7229.5854
$ws.Range("I64").Value = 3949.5
$ws.Range("J64").Value = 7584.189
$ws.Range("K64").Value = 3949.5
$ws.Range("L64").Value = 7584.189
$ws.Range("M64").Value = -3701.5
$ws.Range("N64").Value = -8080.189
$ws.Range("H67").Value = 7229.5854
$ws.Range("I67").Value = 3949.5
$ws.Range("J67").Value = 7584.189
$ws.Range("K67").Value = 3949.5
$ws.Range("L67").Value = 7584.189
$ws.Range("M67").Value = -3091.5
$ws.Range("N67").Value = -9300.189
$ws.Range("H69").Value = 7611.95
$ws.Range("J69").Value = 7846.6113
$ws.Range("L69").Value = 23539.8339
$ws.Range("N69").Value = -25287.8339
$ws.Range("H72").Value = 7611.95
$ws.Range("J72").Value = 7846.6113
$ws.Range("L72").Value = 70619.50169999999
$ws.Range("N72").Value = -79355.50169999999
$ws.Range("H76").Value = 7906.5835
$ws.Range("I76").Value = 7839.857
$ws.Range("K76").Value = 7839.857
$ws.Range("M76").Value = -7524.857
$ws.Range("H79").Value = 7906.5835
$ws.Range("I79").Value = 7839.857
$ws.Range("K79").Value = 7839.857
$ws.Range("M79").Value = -6747.857
$ws.Range("H111").Value = 12350168
$ws.Range("I111").Value = 22226402
$ws.Range("J111").Value = 4875
$ws.Range("K111").Value = 66679206
$ws.Range("L111").Value = 14625
$ws.Range("M111").Value = -66676139
$ws.Range("N111").Value = -20759
$ws.Range("H115").Value = 1971.2858
$ws.Range("I115").Value = 1799.8334
$ws.Range("K115").Value = 5399.5002
$ws.Range("M115").Value = -3832.5002
$ws.Range("H116").Value = 9412.571
$ws.Range("I116").Value = 5000
$ws.Range("J116").Value = 11177.6
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 11177.6
$ws.Range("M116").Value = -1558
$ws.Range("N116").Value = -18061.6
$ws.Range("H133").Value = 92500
$ws.Range("J133").Value = 92500
$ws.Range("L133").Value = 92500
$ws.Range("N133").Value = -102620
$ws.Range("H135").Value = 1580.3334
$ws.Range("I135").Value = 944.6667
$ws.Range("K135").Value = 8502.0003
$ws.Range("M135").Value = -5967.0003
$ws.Range("H137").Value = 63994.45
$ws.Range("I137").Value = 100949.39
$ws.Range("J137").Value = 3522.7273
$ws.Range("K137").Value = 302848.17
$ws.Range("L137").Value = 10568.1819
$ws.Range("M137").Value = -300298.17
$ws.Range("N137").Value = -15668.1819
$ws.Range("H138").Value = 3981.8462
$ws.Range("I138").Value = 3370.5715
$ws.Range("J138").Value = 4115.5625
$ws.Range("K138").Value = 10111.7145
$ws.Range("L138").Value = 12346.6875
$ws.Range("M138").Value = -4971.7145
$ws.Range("N138").Value = -22626.6875
$ws.Range("H141").Value = 2146.9
$ws.Range("I141").Value = 2146.9
$ws.Range("K141").Value = 6440.700000000001
$ws.Range("M141").Value = -1260.700000000001

# ---- ARM (sheet2.xml) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4715690.5
$ws.Range("J2").Value = 2995
$ws.Range("L2").Value = 2995
$ws.Range("N2").Value = -3221
$ws.Range("H32").Value = 3640.6
$ws.Range("I32").Value = 3672.4443
$ws.Range("J32").Value = 488
$ws.Range("K32").Value = 3672.4443
$ws.Range("L32").Value = 488
$ws.Range("M32").Value = -3385.4443
$ws.Range("N32").Value = -1062
$ws.Range("H45").Value = 6807105.5
$ws.Range("I45").Value = 14288093
$ws.Range("J45").Value = 6207.909
$ws.Range("K45").Value = 14288093
$ws.Range("L45").Value = 6207.909
$ws.Range("M45").Value = -14287716
$ws.Range("N45").Value = -6961.909
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H61").Value = 4671.0312
$ws.Range("I61").Value = 5490.077
$ws.Range("J61").Value = 3393.32
$ws.Range("K61").Value = 5490.077
$ws.Range("L61").Value = 3393.32
$ws.Range("M61").Value = -5278.077
$ws.Range("N61").Value = -3817.32
$ws.Range("H63").Value = 6214.769
$ws.Range("I63").Value = 3841.7144
$ws.Range("K63").Value = 3841.7144
$ws.Range("M63").Value = -3155.7144
$ws.Range("H66").Value = 6214.769
$ws.Range("I66").Value = 3841.7144
$ws.Range("K66").Value = 19208.572
$ws.Range("M66").Value = -15776.572
$ws.Range("H74").Value = 21009.068
$ws.Range("I74").Value = 5513.2246
$ws.Range("K74").Value = 5513.2246
$ws.Range("M74").Value = -4639.2246
$ws.Range("H77").Value = 21009.068
$ws.Range("I77").Value = 5513.2246
$ws.Range("K77").Value = 27566.123
$ws.Range("M77").Value = -23198.123
$ws.Range("H88").Value = 2232.6365
$ws.Range("I88").Value = 3331.8
$ws.Range("K88").Value = 3331.8
$ws.Range("M88").Value = -2925.8
$ws.Range("H91").Value = 2232.6365
$ws.Range("I91").Value = 3331.8
$ws.Range("K91").Value = 3331.8
$ws.Range("M91").Value = -1927.8
$ws.Range("H102").Value = 2979261.5
$ws.Range("I102").Value = 3208108.5
$ws.Range("K102").Value = 3208108.5
$ws.Range("M102").Value = -3206486.5
$ws.Range("H110").Value = 993226.2
$ws.Range("I110").Value = 1158459.8
$ws.Range("J110").Value = 1824.5
$ws.Range("K110").Value = 1158459.8
$ws.Range("L110").Value = 1824.5
$ws.Range("M110").Value = -1156414.8
$ws.Range("N110").Value = -5914.5
$ws.Range("H116").Value = 4715690.5
$ws.Range("J116").Value = 2995
$ws.Range("L116").Value = 2995
$ws.Range("N116").Value = -7583
$ws.Range("H122").Value = 3343281
$ws.Range("I122").Value = 4387907
$ws.Range("J122").Value = 2089729.6
$ws.Range("K122").Value = 13163721
$ws.Range("L122").Value = 6269188.800000001
$ws.Range("M122").Value = -13161271
$ws.Range("N122").Value = -6274088.800000001
$ws.Range("H132").Value = 41684.883
$ws.Range("I132").Value = 12039.5
$ws.Range("J132").Value = 60213.25
$ws.Range("K132").Value = 36118.5
$ws.Range("L132").Value = 180639.75
$ws.Range("M132").Value = -33588.5
$ws.Range("N132").Value = -185699.75
$ws.Range("H136").Value = 4671.0312
$ws.Range("I136").Value = 5490.077
$ws.Range("J136").Value = 3393.32
$ws.Range("K136").Value = 16470.231
$ws.Range("L136").Value = 10179.96
$ws.Range("M136").Value = -13920.231
$ws.Range("N136").Value = -15279.96
$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360

# ---- BSM (sheet3.xml) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4715690.5
$ws.Range("J3").Value = 2995
$ws.Range("L3").Value = 2995
$ws.Range("N3").Value = -3223
$ws.Range("H86").Value = 14949065
$ws.Range("I86").Value = 30954598
$ws.Range("J86").Value = 10568.267
$ws.Range("K86").Value = 30954598
$ws.Range("L86").Value = 10568.267
$ws.Range("M86").Value = -30953475
$ws.Range("N86").Value = -12814.267
$ws.Range("H89").Value = 14949065
$ws.Range("I89").Value = 30954598
$ws.Range("J89").Value = 10568.267
$ws.Range("K89").Value = 154772990
$ws.Range("L89").Value = 52841.335
$ws.Range("M89").Value = -154767374
$ws.Range("N89").Value = -64073.335
$ws.Range("H92").Value = 43000
$ws.Range("J92").Value = 43000
$ws.Range("L92").Value = 43000
$ws.Range("N92").Value = -47992
$ws.Range("H105").Value = 3908849.5
$ws.Range("I105").Value = 5211133
$ws.Range("K105").Value = 5211133
$ws.Range("M105").Value = -5209386
$ws.Range("H106").Value = 33000
$ws.Range("J106").Value = 33000
$ws.Range("L106").Value = 33000
$ws.Range("N106").Value = -35524
$ws.Range("H134").Value = 9974.757
$ws.Range("I134").Value = 7783.9062
$ws.Range("J134").Value = 23996.2
$ws.Range("K134").Value = 23351.7186
$ws.Range("L134").Value = 71988.60000000001
$ws.Range("M134").Value = -20816.7186
$ws.Range("N134").Value = -77058.60000000001
$ws.Range("H137").Value = 91978.8
$ws.Range("J137").Value = 92296.25
$ws.Range("L137").Value = 92296.25
$ws.Range("N137").Value = -102496.25
$ws.Range("H140").Value = 67010.28999999999
$ws.Range("J140").Value = 69893.836
$ws.Range("L140").Value = 69893.836
$ws.Range("N140").Value = -80253.836

# ---- CRP (sheet4.xml) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2254.8572
$ws.Range("I16").Value = 1797
$ws.Range("K16").Value = 1797
$ws.Range("M16").Value = -1510
$ws.Range("H31").Value = 29164.5
$ws.Range("I31").Value = 13827.75
$ws.Range("J31").Value = 32998.688
$ws.Range("K31").Value = 13827.75
$ws.Range("L31").Value = 32998.688
$ws.Range("M31").Value = -13532.75
$ws.Range("N31").Value = -33588.688
$ws.Range("H34").Value = 29164.5
$ws.Range("I34").Value = 13827.75
$ws.Range("J34").Value = 32998.688
$ws.Range("K34").Value = 13827.75
$ws.Range("L34").Value = 32998.688
$ws.Range("M34").Value = -13625.75
$ws.Range("N34").Value = -33402.688
$ws.Range("H43").Value = 14477.333
$ws.Range("J43").Value = 14477.333
$ws.Range("L43").Value = 14477.333
$ws.Range("N43").Value = -14845.333
$ws.Range("H58").Value = 8936.9375
$ws.Range("I58").Value = 14358.625
$ws.Range("J58").Value = 3515.25
$ws.Range("K58").Value = 14358.625
$ws.Range("L58").Value = 3515.25
$ws.Range("M58").Value = -14155.625
$ws.Range("N58").Value = -3921.25
$ws.Range("H62").Value = 3023.1667
$ws.Range("I62").Value = 3479.5
$ws.Range("J62").Value = 2110.5
$ws.Range("K62").Value = 3479.5
$ws.Range("L62").Value = 2110.5
$ws.Range("M62").Value = -2855.5
$ws.Range("N62").Value = -3358.5
$ws.Range("H65").Value = 3023.1667
$ws.Range("I65").Value = 3479.5
$ws.Range("J65").Value = 2110.5
$ws.Range("K65").Value = 17397.5
$ws.Range("L65").Value = 10552.5
$ws.Range("M65").Value = -14277.5
$ws.Range("N65").Value = -16792.5
$ws.Range("H86").Value = 5971.231
$ws.Range("I86").Value = 4342.3
$ws.Range("J86").Value = 11401
$ws.Range("K86").Value = 4342.3
$ws.Range("L86").Value = 11401
$ws.Range("M86").Value = -3219.3
$ws.Range("N86").Value = -13647
$ws.Range("H89").Value = 5971.231
$ws.Range("I89").Value = 4342.3
$ws.Range("J89").Value = 11401
$ws.Range("K89").Value = 21711.5
$ws.Range("L89").Value = 57005
$ws.Range("M89").Value = -16095.5
$ws.Range("N89").Value = -68237
$ws.Range("H92").Value = 31666
$ws.Range("J92").Value = 31666
$ws.Range("L92").Value = 31666
$ws.Range("N92").Value = -36658
$ws.Range("H95").Value = 15517.714
$ws.Range("J95").Value = 15517.714
$ws.Range("L95").Value = 15517.714
$ws.Range("N95").Value = -21009.714
$ws.Range("H96").Value = 15312
$ws.Range("J96").Value = 15312
$ws.Range("L96").Value = 15312
$ws.Range("N96").Value = -20804
$ws.Range("H97").Value = 58995
$ws.Range("J97").Value = 58995
$ws.Range("L97").Value = 58995
$ws.Range("N97").Value = -60977
$ws.Range("H101").Value = 14477.333
$ws.Range("J101").Value = 14477.333
$ws.Range("L101").Value = 14477.333
$ws.Range("N101").Value = -20967.333
$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -55242
$ws.Range("H105").Value = 637.6923
$ws.Range("I105").Value = 649.1111
$ws.Range("K105").Value = 649.1111
$ws.Range("M105").Value = 1097.8889
$ws.Range("H107").Value = 1979.2122
$ws.Range("I107").Value = 1882.5714
$ws.Range("K107").Value = 1882.5714
$ws.Range("M107").Value = 37.42859999999996
$ws.Range("H113").Value = 2254.8572
$ws.Range("I113").Value = 1797
$ws.Range("K113").Value = 1797
$ws.Range("M113").Value = 373
$ws.Range("H122").Value = 1940.5714
$ws.Range("I122").Value = 2053.4
$ws.Range("J122").Value = 1658.5
$ws.Range("K122").Value = 6160.200000000001
$ws.Range("L122").Value = 4975.5
$ws.Range("M122").Value = -3710.200000000001
$ws.Range("N122").Value = -9875.5
$ws.Range("H132").Value = 39753.56
$ws.Range("I132").Value = 27932.184
$ws.Range("K132").Value = 83796.552
$ws.Range("M132").Value = -81266.552
$ws.Range("H134").Value = 14254.5
$ws.Range("I134").Value = 13199.111
$ws.Range("J134").Value = 15309.889
$ws.Range("K134").Value = 39597.333
$ws.Range("L134").Value = 45929.667
$ws.Range("M134").Value = -37062.333
$ws.Range("N134").Value = -50999.667
$ws.Range("H136").Value = 8936.9375
$ws.Range("I136").Value = 14358.625
$ws.Range("J136").Value = 3515.25
$ws.Range("K136").Value = 43075.875
$ws.Range("L136").Value = 10545.75
$ws.Range("M136").Value = -40525.875
$ws.Range("N136").Value = -15645.75
$ws.Range("H141").Value = 152299.45
$ws.Range("J141").Value = 166333.22
$ws.Range("L141").Value = 166333.22
$ws.Range("N141").Value = -176693.22

# ---- CUL (sheet5.xml) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2921.0557
$ws.Range("I3").Value = 2229.4
$ws.Range("J3").Value = 3785.625
$ws.Range("K3").Value = 6688.200000000001
$ws.Range("L3").Value = 11356.875
$ws.Range("M3").Value = -6576.200000000001
$ws.Range("N3").Value = -11580.875
$ws.Range("H18").Value = 992.1111
$ws.Range("I18").Value = 238.16667
$ws.Range("J18").Value = 2500
$ws.Range("K18").Value = 714.50001
$ws.Range("L18").Value = 7500
$ws.Range("M18").Value = -545.50001
$ws.Range("N18").Value = -7838
$ws.Range("H29").Value = 359.33334
$ws.Range("I29").Value = 411.2
$ws.Range("K29").Value = 1233.6
$ws.Range("M29").Value = -956.5999999999999
$ws.Range("H55").Value = 62502130
$ws.Range("I55").Value = 124876410
$ws.Range("J55").Value = 127850
$ws.Range("K55").Value = 374629230
$ws.Range("L55").Value = 383550
$ws.Range("M55").Value = -374629053
$ws.Range("N55").Value = -383904
$ws.Range("H68").Value = 1655
$ws.Range("J68").Value = 2125.5
$ws.Range("L68").Value = 6376.5
$ws.Range("N68").Value = -7998.5
$ws.Range("H71").Value = 1655
$ws.Range("J71").Value = 2125.5
$ws.Range("L71").Value = 19129.5
$ws.Range("N71").Value = -27241.5
$ws.Range("H75").Value = 3000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 3000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 9000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -10996
$ws.Range("H78").Value = 3000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 3000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 27000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -36984
$ws.Range("H107").Value = 713.75
$ws.Range("I107").Value = 949
$ws.Range("J107").Value = 478.5
$ws.Range("K107").Value = 2847
$ws.Range("L107").Value = 1435.5
$ws.Range("M107").Value = -927
$ws.Range("N107").Value = -5275.5
$ws.Range("H125").Value = 9300
$ws.Range("I125").Value = 8500
$ws.Range("K125").Value = 25500
$ws.Range("M125").Value = -20580
$ws.Range("H126").Value = 2267.375
$ws.Range("I126").Value = 2187.8
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 6563.400000000001
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -1623.400000000001
$ws.Range("N126").Value = -17080
$ws.Range("H130").Value = 3193.9
$ws.Range("I130").Value = 2848.1667
$ws.Range("J130").Value = 3712.5
$ws.Range("K130").Value = 8544.500100000001
$ws.Range("L130").Value = 11137.5
$ws.Range("M130").Value = -3524.500100000001
$ws.Range("N130").Value = -21177.5
$ws.Range("H132").Value = 1669.9615
$ws.Range("J132").Value = 1657.1177
$ws.Range("L132").Value = 14914.0593
$ws.Range("N132").Value = -19974.0593
$ws.Range("H133").Value = 3467.875
$ws.Range("I133").Value = 3467.875
$ws.Range("K133").Value = 10403.625
$ws.Range("M133").Value = -5343.625
$ws.Range("H134").Value = 2047.1111
$ws.Range("I134").Value = 2047.1111
$ws.Range("K134").Value = 6141.3333
$ws.Range("M134").Value = -1071.3333
$ws.Range("H136").Value = 5783.273
$ws.Range("I136").Value = 5361.6
$ws.Range("K136").Value = 16084.8
$ws.Range("M136").Value = -10984.8
$ws.Range("H137").Value = 4414.1562
$ws.Range("I137").Value = 2368.6155
$ws.Range("J137").Value = 5813.737
$ws.Range("K137").Value = 7105.8465
$ws.Range("L137").Value = 17441.211
$ws.Range("M137").Value = -2005.8465
$ws.Range("N137").Value = -27641.211
$ws.Range("H138").Value = 4229.8887
$ws.Range("I138").Value = 2834
$ws.Range("J138").Value = 5974.75
$ws.Range("K138").Value = 8502
$ws.Range("L138").Value = 17924.25
$ws.Range("M138").Value = -3362
$ws.Range("N138").Value = -28204.25
$ws.Range("H139").Value = 2787.4285
$ws.Range("I139").Value = 2787.4285
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 8362.2855
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -3222.2855
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 1628
$ws.Range("I140").Value = 1035
$ws.Range("K140").Value = 3105
$ws.Range("M140").Value = 2075

# ---- GSM (sheet6.xml) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10536398
$ws.Range("I70").Value = 18185808
$ws.Range("K70").Value = 18185808
$ws.Range("M70").Value = -18185538
$ws.Range("H73").Value = 10536398
$ws.Range("I73").Value = 18185808
$ws.Range("K73").Value = 18185808
$ws.Range("M73").Value = -18184872
$ws.Range("H80").Value = 1820284.5
$ws.Range("I80").Value = 2750474
$ws.Range("J80").Value = 425000
$ws.Range("K80").Value = 2750474
$ws.Range("L80").Value = 425000
$ws.Range("M80").Value = -2749476
$ws.Range("N80").Value = -426996
$ws.Range("H83").Value = 1820284.5
$ws.Range("I83").Value = 2750474
$ws.Range("J83").Value = 425000
$ws.Range("K83").Value = 13752370
$ws.Range("L83").Value = 2125000
$ws.Range("M83").Value = -13747378
$ws.Range("N83").Value = -2134984
$ws.Range("H97").Value = 916528.4
$ws.Range("I97").Value = 1323406.1
$ws.Range("J97").Value = 1053.375
$ws.Range("K97").Value = 1323406.1
$ws.Range("L97").Value = 1053.375
$ws.Range("M97").Value = -1322910.1
$ws.Range("N97").Value = -2045.375
$ws.Range("H102").Value = 7498213.5
$ws.Range("I102").Value = 13890192
$ws.Range("K102").Value = 13890192
$ws.Range("M102").Value = -13888570
$ws.Range("H107").Value = 10729.5
$ws.Range("I107").Value = 17284.166
$ws.Range("K107").Value = 17284.166
$ws.Range("M107").Value = -15364.166
$ws.Range("H113").Value = 16064179
$ws.Range("I113").Value = 35335616
$ws.Range("J113").Value = 4650
$ws.Range("K113").Value = 35335616
$ws.Range("L113").Value = 4650
$ws.Range("M113").Value = -35333446
$ws.Range("N113").Value = -8990
$ws.Range("H122").Value = 1114273.2
$ws.Range("I122").Value = 1485366
$ws.Range("J122").Value = 995
$ws.Range("K122").Value = 4456098
$ws.Range("L122").Value = 2985
$ws.Range("M122").Value = -4453648
$ws.Range("N122").Value = -7885
$ws.Range("H132").Value = 6391.8394
$ws.Range("I132").Value = 4998.1904
$ws.Range("J132").Value = 10572.786
$ws.Range("K132").Value = 14994.5712
$ws.Range("L132").Value = 31718.358
$ws.Range("M132").Value = -12464.5712
$ws.Range("N132").Value = -36778.358
$ws.Range("H138").Value = 73899.8
$ws.Range("J138").Value = 73899.8
$ws.Range("L138").Value = 73899.8
$ws.Range("N138").Value = -84179.8
$ws.Range("H140").Value = 105866.445
$ws.Range("J140").Value = 105866.445
$ws.Range("L140").Value = 105866.445
$ws.Range("N140").Value = -116226.445

# ---- LTW (sheet7.xml) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2895476.2
$ws.Range("J2").Value = 53666.8
$ws.Range("L2").Value = 53666.8
$ws.Range("N2").Value = -53890.8
$ws.Range("H22").Value = 150583.5
$ws.Range("J22").Value = 4334
$ws.Range("L22").Value = 4334
$ws.Range("N22").Value = -4924
$ws.Range("H27").Value = 150583.5
$ws.Range("J27").Value = 4334
$ws.Range("L27").Value = 4334
$ws.Range("N27").Value = -4548
$ws.Range("H40").Value = 6008.2646
$ws.Range("J40").Value = 12298.2
$ws.Range("L40").Value = 12298.2
$ws.Range("N40").Value = -12570.2
$ws.Range("H46").Value = 4977.875
$ws.Range("I46").Value = 1350.5
$ws.Range("J46").Value = 5307.636
$ws.Range("K46").Value = 1350.5
$ws.Range("L46").Value = 5307.636
$ws.Range("M46").Value = -1162.5
$ws.Range("N46").Value = -5683.636
$ws.Range("H61").Value = 10101813
$ws.Range("I61").Value = 10101813
$ws.Range("K61").Value = 10101813
$ws.Range("M61").Value = -10101611
$ws.Range("H68").Value = 1999.125
$ws.Range("I68").Value = 1638
$ws.Range("J68").Value = 2601
$ws.Range("K68").Value = 1638
$ws.Range("L68").Value = 2601
$ws.Range("M68").Value = -889
$ws.Range("N68").Value = -4099
$ws.Range("H71").Value = 1999.125
$ws.Range("I71").Value = 1638
$ws.Range("J71").Value = 2601
$ws.Range("K71").Value = 8190
$ws.Range("L71").Value = 13005
$ws.Range("M71").Value = -4446
$ws.Range("N71").Value = -20493
$ws.Range("H100").Value = 3326.8965
$ws.Range("I100").Value = 2712.8572
$ws.Range("J100").Value = 3900
$ws.Range("K100").Value = 2712.8572
$ws.Range("L100").Value = 3900
$ws.Range("M100").Value = -2171.8572
$ws.Range("N100").Value = -4982
$ws.Range("H113").Value = 10101813
$ws.Range("I113").Value = 10101813
$ws.Range("K113").Value = 10101813
$ws.Range("M113").Value = -10099643
$ws.Range("H122").Value = 7864.6665
$ws.Range("I122").Value = 4996.6665
$ws.Range("J122").Value = 9298.666999999999
$ws.Range("K122").Value = 14989.9995
$ws.Range("L122").Value = 27896.001
$ws.Range("M122").Value = -12539.9995
$ws.Range("N122").Value = -32796.001
$ws.Range("H130").Value = 55212
$ws.Range("J130").Value = 55212
$ws.Range("L130").Value = 55212
$ws.Range("N130").Value = -65252
$ws.Range("H132").Value = 14620.679
$ws.Range("I132").Value = 15437.654
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 46312.962
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -43782.962
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 81281.19500000001
$ws.Range("I136").Value = 108437.42
$ws.Range("K136").Value = 325312.26
$ws.Range("M136").Value = -322762.26
$ws.Range("H137").Value = 140000
$ws.Range("J137").Value = 140000
$ws.Range("L137").Value = 140000
$ws.Range("N137").Value = -150200
$ws.Range("H141").Value = 88552.86
$ws.Range("J141").Value = 88552.86
$ws.Range("L141").Value = 88552.86
$ws.Range("N141").Value = -98912.86

# ---- WVR (sheet8.xml) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 19995
$ws.Range("J32").Value = 19995
$ws.Range("L32").Value = 19995
$ws.Range("N32").Value = -20629
$ws.Range("H51").Value = 10023.333
$ws.Range("J51").Value = 10000
$ws.Range("L51").Value = 10000
$ws.Range("N51").Value = -11020
$ws.Range("H62").Value = 11034.272
$ws.Range("I62").Value = 16658.666
$ws.Range("J62").Value = 8925.125
$ws.Range("K62").Value = 16658.666
$ws.Range("L62").Value = 8925.125
$ws.Range("M62").Value = -16034.666
$ws.Range("N62").Value = -10173.125
$ws.Range("H65").Value = 11034.272
$ws.Range("I65").Value = 16658.666
$ws.Range("J65").Value = 8925.125
$ws.Range("K65").Value = 83293.33
$ws.Range("L65").Value = 44625.625
$ws.Range("M65").Value = -80173.33
$ws.Range("N65").Value = -50865.625
$ws.Range("H107").Value = 41667580
$ws.Range("I107").Value = 71429130
$ws.Range("K107").Value = 214287390
$ws.Range("M107").Value = -214285470
$ws.Range("H132").Value = 19082726
$ws.Range("I132").Value = 22230584
$ws.Range("K132").Value = 66691752
$ws.Range("M132").Value = -66689222
$ws.Range("H133").Value = 59388.715
$ws.Range("J133").Value = 59388.715
$ws.Range("L133").Value = 59388.715
$ws.Range("N133").Value = -69508.715
$ws.Range("H136").Value = 5495.476
$ws.Range("I136").Value = 5921.5757
$ws.Range("J136").Value = 3933.111
$ws.Range("K136").Value = 17764.7271
$ws.Range("L136").Value = 11799.333
$ws.Range("M136").Value = -15214.7271
$ws.Range("N136").Value = -16899.333

